# Weekly update: insert a new daily price record for Perejil (Feria Lagunitas
# de Puerto Montt) on 2022-02-09 (serial 44601), pushing the existing
# historical rows down by one (row 98 -> 99, ..., old row 209 -> 210).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right above the current row 98; this shifts every
# row from 98..209 down to 99..210 and grows the used range to A1:R210.
$ws.Rows(98).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A98").Value = 4
$ws.Range("B98").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C98").Value = "Los Lagos"
$ws.Range("D98").Value = 44601
$ws.Range("D98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E98").Value = 10
$ws.Range("F98").Value = 100112044
$ws.Range("G98").Value = "Perejil"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 40
$ws.Range("K98").Value = 6000
$ws.Range("L98").Value = 6000
$ws.Range("M98").Value = 6000
$ws.Range("N98").Value = "$/docena de atados (3 kilos)"
$ws.Range("O98").Value = "Región Metropolitana"
$ws.Range("P98").Value = 2000
$ws.Range("Q98").Value = 3
$ws.Range("R98").Value = "Hortaliza"
